$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Translate the Spanish destination names that remain in the sheet to
# --- their English equivalents (Praga/Copenhague/Estambul/Bruselas/
# --- Edimburgo/Estocolmo/Lisboa/Roma -> Prague/Copenhagen/Istanbul/
# --- Brussels/Edinburgh/Stockholm/Lisbon/Rome). Every cell below held one
# --- of those Spanish city names.
$translations = @{
    "B11"  = "Rome"
    "C15"  = "Brussels"
    "B16"  = "Brussels"
    "C17"  = "Brussels"
    "B18"  = "Brussels"
    "C19"  = "Brussels"
    "B22"  = "Brussels"
    "C120" = "Stockholm"
    "B121" = "Stockholm"
    "C137" = "Prague"
    "B138" = "Prague"
    "C141" = "Lisbon"
    "B142" = "Lisbon"
    "C146" = "Copenhagen"
    "B147" = "Copenhagen"
    "C148" = "Istanbul"
    "B149" = "Istanbul"
    "C150" = "Istanbul"
    "B151" = "Istanbul"
    "C153" = "Brussels"
    "B154" = "Brussels"
    "C155" = "Edinburgh"
    "B156" = "Edinburgh"
    "C157" = "Edinburgh"
    "B158" = "Edinburgh"
    "C159" = "Edinburgh"
    "B160" = "Edinburgh"
    "C161" = "Edinburgh"
    "B162" = "Edinburgh"
    "C163" = "Edinburgh"
    "B164" = "Edinburgh"
    "C165" = "Edinburgh"
    "B166" = "Edinburgh"
    "B168" = "Edinburgh"
}

foreach ($addr in $translations.Keys) {
    $ws.Range($addr).Value = $translations[$addr]
}

# --- Append a new round-trip flight record (id 88): Madrid -> Palma de
# --- Mallorca -> Madrid, August 2019, Holiday.
$ws.Range("A185").Formula = "=A183+1"
$ws.Range("B185").Value = "Madrid"
$ws.Range("C185").Value = "Palma de Mallorca"
$ws.Range("D185").Value = 8
$ws.Range("E185").Value = 2019
$ws.Range("F185").Value = "Holiday"

$ws.Range("A186").Formula = "=A184+1"
$ws.Range("B186").Value = "Palma de Mallorca"
$ws.Range("C186").Value = "Madrid"
$ws.Range("D186").Value = 8
$ws.Range("E186").Value = 2019
$ws.Range("F186").Value = "Holiday"

# --- Update the saved view state: zoom to 85%, drop the old frozen
# --- top-left cell, and move the active selection to K10.
$ws.Activate()
$excel.ActiveWindow.Zoom = 85
$ws.Range("K10").Select()
